$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

if ($t.Cell(1,1).Range.Text.TrimEnd([char]7,[char]13) -ne "43÷8=5, 3") { throw "Mismatch at (1,1): expected '43÷8=5, 3' but found '$($t.Cell(1,1).Range.Text)'" }
$t.Cell(1,1).Range.Text = "86÷2=43, 0"

if ($t.Cell(1,2).Range.Text.TrimEnd([char]7,[char]13) -ne "64÷8=8, 0") { throw "Mismatch at (1,2): expected '64÷8=8, 0' but found '$($t.Cell(1,2).Range.Text)'" }
$t.Cell(1,2).Range.Text = "17÷7=2, 3"

if ($t.Cell(1,3).Range.Text.TrimEnd([char]7,[char]13) -ne "89÷8=11, 1") { throw "Mismatch at (1,3): expected '89÷8=11, 1' but found '$($t.Cell(1,3).Range.Text)'" }
$t.Cell(1,3).Range.Text = "88÷4=22, 0"

if ($t.Cell(1,4).Range.Text.TrimEnd([char]7,[char]13) -ne "72÷5=14, 2") { throw "Mismatch at (1,4): expected '72÷5=14, 2' but found '$($t.Cell(1,4).Range.Text)'" }
$t.Cell(1,4).Range.Text = "45÷4=11, 1"

if ($t.Cell(1,5).Range.Text.TrimEnd([char]7,[char]13) -ne "54÷8=6, 6") { throw "Mismatch at (1,5): expected '54÷8=6, 6' but found '$($t.Cell(1,5).Range.Text)'" }
$t.Cell(1,5).Range.Text = "28÷5=5, 3"

if ($t.Cell(5,1).Range.Text.TrimEnd([char]7,[char]13) -ne "26÷3=8, 2") { throw "Mismatch at (5,1): expected '26÷3=8, 2' but found '$($t.Cell(5,1).Range.Text)'" }
$t.Cell(5,1).Range.Text = "52÷8=6, 4"

if ($t.Cell(5,2).Range.Text.TrimEnd([char]7,[char]13) -ne "65÷2=32, 1") { throw "Mismatch at (5,2): expected '65÷2=32, 1' but found '$($t.Cell(5,2).Range.Text)'" }
$t.Cell(5,2).Range.Text = "95÷2=47, 1"

if ($t.Cell(5,3).Range.Text.TrimEnd([char]7,[char]13) -ne "93÷6=15, 3") { throw "Mismatch at (5,3): expected '93÷6=15, 3' but found '$($t.Cell(5,3).Range.Text)'" }
$t.Cell(5,3).Range.Text = "73÷3=24, 1"

if ($t.Cell(5,4).Range.Text.TrimEnd([char]7,[char]13) -ne "89÷4=22, 1") { throw "Mismatch at (5,4): expected '89÷4=22, 1' but found '$($t.Cell(5,4).Range.Text)'" }
$t.Cell(5,4).Range.Text = "46÷6=7, 4"

if ($t.Cell(5,5).Range.Text.TrimEnd([char]7,[char]13) -ne "25÷9=2, 7") { throw "Mismatch at (5,5): expected '25÷9=2, 7' but found '$($t.Cell(5,5).Range.Text)'" }
$t.Cell(5,5).Range.Text = "56÷6=9, 2"

if ($t.Cell(9,1).Range.Text.TrimEnd([char]7,[char]13) -ne "75÷2=37, 1") { throw "Mismatch at (9,1): expected '75÷2=37, 1' but found '$($t.Cell(9,1).Range.Text)'" }
$t.Cell(9,1).Range.Text = "22÷8=2, 6"

if ($t.Cell(9,2).Range.Text.TrimEnd([char]7,[char]13) -ne "37÷5=7, 2") { throw "Mismatch at (9,2): expected '37÷5=7, 2' but found '$($t.Cell(9,2).Range.Text)'" }
$t.Cell(9,2).Range.Text = "36÷5=7, 1"

if ($t.Cell(9,3).Range.Text.TrimEnd([char]7,[char]13) -ne "37÷8=4, 5") { throw "Mismatch at (9,3): expected '37÷8=4, 5' but found '$($t.Cell(9,3).Range.Text)'" }
$t.Cell(9,3).Range.Text = "34÷5=6, 4"

if ($t.Cell(9,4).Range.Text.TrimEnd([char]7,[char]13) -ne "67÷2=33, 1") { throw "Mismatch at (9,4): expected '67÷2=33, 1' but found '$($t.Cell(9,4).Range.Text)'" }
$t.Cell(9,4).Range.Text = "59÷3=19, 2"

if ($t.Cell(9,5).Range.Text.TrimEnd([char]7,[char]13) -ne "15÷5=3, 0") { throw "Mismatch at (9,5): expected '15÷5=3, 0' but found '$($t.Cell(9,5).Range.Text)'" }
$t.Cell(9,5).Range.Text = "49÷4=12, 1"

if ($t.Cell(13,1).Range.Text.TrimEnd([char]7,[char]13) -ne "85÷4=21, 1") { throw "Mismatch at (13,1): expected '85÷4=21, 1' but found '$($t.Cell(13,1).Range.Text)'" }
$t.Cell(13,1).Range.Text = "55÷2=27, 1"

if ($t.Cell(13,2).Range.Text.TrimEnd([char]7,[char]13) -ne "24÷6=4, 0") { throw "Mismatch at (13,2): expected '24÷6=4, 0' but found '$($t.Cell(13,2).Range.Text)'" }
$t.Cell(13,2).Range.Text = "92÷6=15, 2"

if ($t.Cell(13,3).Range.Text.TrimEnd([char]7,[char]13) -ne "24÷9=2, 6") { throw "Mismatch at (13,3): expected '24÷9=2, 6' but found '$($t.Cell(13,3).Range.Text)'" }
$t.Cell(13,3).Range.Text = "71÷9=7, 8"

if ($t.Cell(13,4).Range.Text.TrimEnd([char]7,[char]13) -ne "26÷6=4, 2") { throw "Mismatch at (13,4): expected '26÷6=4, 2' but found '$($t.Cell(13,4).Range.Text)'" }
$t.Cell(13,4).Range.Text = "38÷4=9, 2"

if ($t.Cell(13,5).Range.Text.TrimEnd([char]7,[char]13) -ne "28÷4=7, 0") { throw "Mismatch at (13,5): expected '28÷4=7, 0' but found '$($t.Cell(13,5).Range.Text)'" }
$t.Cell(13,5).Range.Text = "62÷7=8, 6"

if ($t.Cell(17,1).Range.Text.TrimEnd([char]7,[char]13) -ne "89÷4=22, 1") { throw "Mismatch at (17,1): expected '89÷4=22, 1' but found '$($t.Cell(17,1).Range.Text)'" }
$t.Cell(17,1).Range.Text = "98÷2=49, 0"

if ($t.Cell(17,2).Range.Text.TrimEnd([char]7,[char]13) -ne "85÷8=10, 5") { throw "Mismatch at (17,2): expected '85÷8=10, 5' but found '$($t.Cell(17,2).Range.Text)'" }
$t.Cell(17,2).Range.Text = "48÷4=12, 0"

if ($t.Cell(17,3).Range.Text.TrimEnd([char]7,[char]13) -ne "75÷5=15, 0") { throw "Mismatch at (17,3): expected '75÷5=15, 0' but found '$($t.Cell(17,3).Range.Text)'" }
$t.Cell(17,3).Range.Text = "20÷5=4, 0"

if ($t.Cell(17,4).Range.Text.TrimEnd([char]7,[char]13) -ne "63÷7=9, 0") { throw "Mismatch at (17,4): expected '63÷7=9, 0' but found '$($t.Cell(17,4).Range.Text)'" }
$t.Cell(17,4).Range.Text = "33÷9=3, 6"

if ($t.Cell(17,5).Range.Text.TrimEnd([char]7,[char]13) -ne "80÷2=40, 0") { throw "Mismatch at (17,5): expected '80÷2=40, 0' but found '$($t.Cell(17,5).Range.Text)'" }
$t.Cell(17,5).Range.Text = "96÷3=32, 0"
